$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell address -> refreshed text value (coinranking.com snapshot refresh).
# Price/volume columns are stored as plain text in this sheet (e.g.
# "2.617.27", "  -0.21%  "), so every write goes through as a string.
# A handful of new Price values happen to look like plain decimals
# (e.g. "7.40", "0.999") -- for those we briefly force a Text number
# format so Excel doesn't silently coerce them into numbers (which would
# drop the trailing zero / change "1.50" into 1.5), then clear the
# temporary format back off so the cell's style is untouched.
$updates = [ordered]@{
    'D2' = '60.635.53'
    'E2' = '  -0.36%  '
    'D3' = '2.617.27'
    'E3' = '  -0.21%  '
    'E4' = '  -0.09%  '
    'D5' = '511.91'
    'E5' = '  +0.36%  '
    'D6' = '155.14'
    'E6' = '  -1.66%  '
    'D7' = '0.997'
    'E7' = '  +0.07%  '
    'E8' = '  -2.90%  '
    'D9' = '2.632.97'
    'E9' = '  -1.05%  '
    'E10' = '  +4.75%  '
    'D11' = '0.105'
    'E11' = '  -0.30%  '
    'D12' = '0.347'
    'E12' = '  -0.16%  '
    'E13' = '  +1.52%  '
    'D14' = '3.075.07'
    'E14' = '  -0.63%  '
    'D15' = '60.571.64'
    'E15' = '  -0.40%  '
    'D16' = '21.66'
    'E16' = '  -0.92%  '
    'E17' = '  -0.10%  '
    'D18' = '2.627.70'
    'E18' = '  -0.94%  '
    'D20' = '351.29'
    'E20' = '  +0.70%  '
    'D21' = '10.63'
    'E21' = '  +0.87%  '
    'E22' = '  -0.33%  '
    'E23' = '  +0.11%  '
    'E24' = '  +0.04%  '
    'D25' = '0.424'
    'E25' = '  -0.07%  '
    'D26' = '0.167'
    'E26' = '  -0.48%  '
    'E27' = '  +1.09%  '
    'D28' = '0.0₃0846'
    'E28' = '  -2.80%  '
    'D29' = '7.40'
    'E29' = '  -2.13%  '
    'D30' = '0.999'
    'E30' = '  +0.15%  '
    'D31' = '19.47'
    'E31' = '  -0.41%  '
    'D32' = '1.58'
    'E32' = '  +0.18%  '
    'D33' = '150.29'
    'E33' = '  -4.16%  '
    'E34' = '  +0.96%  '
    'E35' = '  -1.83%  '
    'E36' = '  -2.03%  '
    'D37' = '0.894'
    'E37' = '  +5.69%  '
    'D38' = '1.50'
    'E38' = '  +0.01%  '
    'E39' = '  -1.00%  '
    'B40' = 'Filecoin'
    'C40' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D40' = '3.78'
    'E40' = '  +0.26%  '
    'B41' = 'OKB'
    'C41' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D41' = '36.40'
    'E41' = '  +3.43%  '
    'D42' = '294.35'
    'E42' = '  -5.76%  '
    'D43' = '0.627'
    'E43' = '  -2.66%  '
    'E44' = '  -0.06%  '
    'D45' = '0.997'
    'E45' = '  +0.09%  '
    'D46' = '19.93'
    'E46' = '  -1.26%  '
    'E47' = '  -4.63%  '
    'D48' = '4.92'
    'E48' = '  +0.66%  '
    'E49' = '  -0.81%  '
    'E50' = '  +0.35%  '
    'D51' = '2.006.23'
    'E51' = '  -2.79%  '
}

foreach ($addr in $updates.Keys) {
    $newValue = $updates[$addr]
    $cell = $ws.Range($addr)

    if ($newValue -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.ClearFormats()
    } else {
        $cell.Value = $newValue
    }
}
